$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2026-02-20 Friday" "2026-02-21 Saturday"

Replace-Text "374÷8=" "846÷9="
Replace-Text "251÷8=" "115÷3="
Replace-Text "643÷6=" "319÷4="
Replace-Text "844÷8=" "705÷6="
Replace-Text "302÷9=" "993÷6="
Replace-Text "812÷7=" "183÷3="
Replace-Text "232÷6=" "989÷5="
Replace-Text "852÷6=" "832÷9="
Replace-Text "763÷5=" "837÷4="
Replace-Text "871÷4=" "288÷3="
Replace-Text "135÷6=" "644÷8="
Replace-Text "252÷2=" "655÷5="
Replace-Text "697÷5=" "148÷2="
Replace-Text "660÷6=" "142÷2="
Replace-Text "245÷2=" "800÷5="
Replace-Text "834÷7=" "819÷9="
Replace-Text "905÷6=" "944÷4="
Replace-Text "373÷8=" "556÷4="
Replace-Text "392÷7=" "963÷7="
Replace-Text "307÷5=" "262÷5="
Replace-Text "625÷7=" "881÷8="
Replace-Text "153÷2=" "109÷2="
Replace-Text "167÷8=" "711÷5="
Replace-Text "131÷5=" "576÷7="
Replace-Text "293÷4=" "846÷7="

Write-Output "Done"
